$p = $ppt.ActivePresentation

# Slide 5 ("Advantages of Application :") holds the bulleted list we need
# to edit; find it by its title text so the script is resilient to any
# slide re-ordering.
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    $title = ""
    try { $title = $candidate.Shapes.Item(1).TextFrame.TextRange.Text } catch {}
    if ($title -like "Advantages of Application*") {
        $s = $candidate
        break
    }
}
if ($s -eq $null) { $s = $p.Slides.Item(5) }

$shp = $s.Shapes.Item(2)
$tf  = $shp.TextFrame
$tr  = $tf.TextRange

# --- Rebuild the body text with the two new bullets added and the
#     trailing "Makes your business approachable" bullet removed. ---
$CR = [char]13

$bullets = @(
    "Real time text preview.",
    "User can not chat without being friend.",
    "Easily share images and others supported files.",
    "User can mute any chat and moved into already created mute folder.",
    "User can also chat in group.",
    "It is fast than email for users. ",
    "Data security.",
    "Low cost.",
    "Easy to monitor."
)

$tr.Text = [string]::Join($CR, $bullets)

# --- Turn on "Shrink text on overflow" (adds <a:normAutofit/>) ---
$tf.AutoSize = 2

# --- Split the "mute" bullet into its three original runs ---
$part1 = "User can mute any chat and moved into already created "
$part2 = "mute "
$part3 = "folder."

$prefixLen = 0
foreach ($b in $bullets[0..2]) { $prefixLen += $b.Length + 1 }

$r1 = $tr.Characters($prefixLen + 1, $part1.Length)
$r1.Text = $part1

$r2 = $tr.Characters($prefixLen + 1 + $part1.Length, $part2.Length)
$r2.Text = $part2

$r3 = $tr.Characters($prefixLen + 1 + $part1.Length + $part2.Length, $part3.Length)
$r3.Text = $part3
